$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22
$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# "25" looks numeric, so force the cell to Text format first - otherwise
# Excel auto-converts it to a Number, unlike the source data (stored as a
# string, matching the other rows in this column).
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
